# Simplify the GST report: rename the sheet, drop the frozen header pane,
# delete the TOTAL / footer rows, and strip all per-cell styling (fonts,
# fills, borders, custom row heights) so the sheet is plain structured data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename "GST Audit Report" -> "GST Report"
$ws.Name = "GST Report"

# 2) Remove the frozen header row / split pane, leaving a plain selection.
$excel.ActiveWindow.FreezePanes = $false

# 3) Drop the TOTAL row (3), the blank row (4) and the footer rows (5-6) —
#    this also removes the A3:E3 merged cell along with row 3.
$ws.Rows("3:6").Delete()

# 4) Strip all styling (fonts/fills/borders/alignment) from the remaining
#    header + data rows, reverting every cell back to the workbook's
#    default "Normal" style.
$ws.Range("A1:K2").Style = "Normal"

# 5) Drop the custom header/data row heights so rows fall back to the
#    sheet's default row height.
$ws.Rows("1:2").AutoFit()

Write-Host "GST report simplified: renamed sheet, removed styling/rows"
